$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting the old "Universal/Version 3d"
# column (and the copied RESPONSE values) one column to the right, to make
# room for the new "CANACE8C" module results.
$ws.Columns("D:D").Insert()

# New column D header: module name + version
$ws.Range("D1").Value = "CANACE8C"
$ws.Range("D2").Value = "Version 2q"

# Copy the RESPONSE 1 / RESPONSE 2 values across into the new column, same as
# the existing "CANACC5" column C.
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("D6").Value = $ws.Range("C6").Value2

# Match column widths from the diff as closely as possible: C and D share
# the width that C had before (E keeps the old D width automatically since
# it shifted with the insert).
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth

# Update the stored selection to match the authored workbook state.
$ws.Range("D12").Select()
